# "Categories and features are updated"
#
# Adds a new feature row (#33 "Geometric mean of search result frequency and
# max Lucene score") to the features sheet, widens column B to fit the new
# (longer) description text, and gives the whole used range a thin box
# border with the header row (row 1) in bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the new feature row ------------------------------------------------
$lastRow = 34
$ws.Range("A$lastRow").Value = 33
$ws.Range("B$lastRow").Value = "Geometric mean of search result frequency and max Lucene score"

# --- 2. Widen column B so the longer description still fits -----------------------
$ws.Columns.Item(2).ColumnWidth = 60.5703125

# --- 3. Box border around every cell in the (now 34-row) table --------------------
$dataRange = $ws.Range("A1:B$lastRow")
$dataRange.Borders.LineStyle = 1
$dataRange.Borders.Weight = 2

# --- 4. Bold the header row --------------------------------------------------------
$ws.Range("A1:B1").Font.Bold = $true

# --- 5. Move the selection to the newly added cell, matching the saved view -------
$ws.Range("B$lastRow").Select() | Out-Null

Write-Output "Added row $lastRow and restyled A1:B$lastRow"
